$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3086
$ws.Range("F3").Value = 495
$ws.Range("F4").Value = 63
$ws.Range("F6").Value = 9
$ws.Range("F7").Value = 266
$ws.Range("F8").Value = 16
$ws.Range("F9").Value = 1068
$ws.Range("F10").Value = 15039
$ws.Range("F11").Value = 192
$ws.Range("F14").Value = 5969
$ws.Range("F15").Value = 608
$ws.Range("F16").Value = 89
$ws.Range("F18").Value = 94
$ws.Range("F21").Value = 100
$ws.Range("F22").Value = 4
$ws.Range("F23").Value = 202
$ws.Range("F24").Value = 828
$ws.Range("F25").Value = 2962
$ws.Range("F27").Value = 10798
$ws.Range("F28").Value = 1214
$ws.Range("F29").Value = 88
$ws.Range("F30").Value = 132
$ws.Range("F31").Value = 3762

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3086
$ws.Range("F4").Value = 495
$ws.Range("F5").Value = 63
$ws.Range("F7").Value = 9
$ws.Range("F8").Value = 266
$ws.Range("F9").Value = 16
$ws.Range("F10").Value = 1068
$ws.Range("F11").Value = 15039
$ws.Range("F12").Value = 192
$ws.Range("F15").Value = 5969
$ws.Range("F16").Value = 608
$ws.Range("F17").Value = 89
$ws.Range("F19").Value = 94
$ws.Range("F22").Value = 100
$ws.Range("F23").Value = 4
$ws.Range("F24").Value = 202
$ws.Range("F25").Value = 828
$ws.Range("F26").Value = 0
$ws.Range("F29").Value = 10798
$ws.Range("F30").Value = 1214
$ws.Range("F31").Value = 88
$ws.Range("F32").Value = 132
$ws.Range("F33").Value = 3762
